$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    "a liberal",
    "a conservative",
    "my mom",
    "your mom",
    "covfefe",
    "a can of cold spaghettiOs",
    "cracking open a cold one with the boys",
    "Super Mario 64",
    "Breath of the Wild",
    "Princess Peach",
    "Luigi",
    "Final Fantasy VIII",
    "Final Fantasy VII",
    "Final Fantasy IX",
    "Donkey Kong",
    "Ocarina of Time",
    "A Link to the Past",
    "Gamecube",
    "Pokemon",
    "that picture you posted on /gonewild",
    "your face",
    "the Jolly Rancher story",
    "the Doritos story",
    "that story about the broken arms",
    "the gay agenda"
)

$startRow = 49
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}

$lastRow = $startRow + $newValues.Length
$ws.Cells.Item($lastRow, 1).Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 44
